# Apply cryptos list update (price/volume refresh + two row content swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.151.30'
$ws.Range("D3").Value = '2.877.69'
$ws.Range("E3").Value = '  -5.68%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '486.77'
$c.ClearFormats()
$ws.Range("E5").Value = '  -6.11%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '130.84'
$c.ClearFormats()
$ws.Range("E6").Value = '  -7.78%  '
$ws.Range("E7").Value = '  +0.04%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.414'
$c.ClearFormats()
$ws.Range("E8").Value = '  -6.85%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.ClearFormats()
$ws.Range("E9").Value = '  -5.62%  '
$ws.Range("E10").Value = '  -8.37%  '
$ws.Range("E11").Value = '  -7.20%  '
$ws.Range("D12").Value = '3.373.15'
$ws.Range("E12").Value = '  -5.68%  '
$ws.Range("E13").Value = '  -4.47%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.20'
$c.ClearFormats()
$ws.Range("E14").Value = '  -5.76%  '
$ws.Range("E15").Value = '  -8.54%  '
$ws.Range("D16").Value = '55.133.79'
$ws.Range("E16").Value = '  -4.94%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.92'
$c.ClearFormats()
$ws.Range("E17").Value = '  -5.18%  '
$ws.Range("D18").Value = '2.875.36'
$ws.Range("E18").Value = '  -5.81%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.27'
$c.ClearFormats()
$ws.Range("E19").Value = '  -5.91%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.53'
$c.ClearFormats()
$ws.Range("E20").Value = '  -6.96%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '310.22'
$c.ClearFormats()
$ws.Range("E21").Value = '  -8.32%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.44%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.08%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.474'
$c.ClearFormats()
$ws.Range("E24").Value = '  -5.56%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '61.65'
$c.ClearFormats()
$ws.Range("E25").Value = '  -5.30%  '
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  -6.80%  '
$ws.Range("D28").Value = '0.0₃0825'
$ws.Range("E28").Value = '  -13.62%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.24'
$c.ClearFormats()
$ws.Range("E29").Value = '  -9.74%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.86'
$c.ClearFormats()
$ws.Range("E30").Value = '  -9.53%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.71'
$c.ClearFormats()
$ws.Range("E31").Value = '  -6.32%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '19.44'
$c.ClearFormats()
$ws.Range("E32").Value = '  -7.62%  '
$ws.Range("E33").Value = '  -10.89%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '147.99'
$c.ClearFormats()
$ws.Range("E34").Value = '  -5.39%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.31'
$c.ClearFormats()
$ws.Range("E35").Value = '  -9.67%  '
$ws.Range("E36").Value = '  -7.42%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '24.13'
$c.ClearFormats()
$ws.Range("E37").Value = '  -4.56%  '
$ws.Range("E38").Value = '  -10.52%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0642'
$c.ClearFormats()
$ws.Range("E39").Value = '  -7.29%  '
$ws.Range("E40").Value = '  -0.07%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '35.91'
$c.ClearFormats()
$ws.Range("E41").Value = '  -4.98%  '
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.31'
$c.ClearFormats()
$ws.Range("E44").Value = '  -9.90%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.066.81'
$ws.Range("E45").Value = '  -11.32%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '5.78'
$c.ClearFormats()
$ws.Range("E46").Value = '  -4.94%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.896'
$c.ClearFormats()
$ws.Range("E47").Value = '  -11.26%  '
$ws.Range("E48").Value = '  -6.23%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '18.34'
$c.ClearFormats()
$ws.Range("E49").Value = '  -7.50%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0834'
$c.ClearFormats()
$ws.Range("E50").Value = '  -7.35%  '
$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '4.93'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.52%  '
